$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: append new trailing data (day-over-day price change + Up/Down verdict)
$ws.Range("X3").Value = -0.93999500000001035
$ws.Range("Y3").Value = "Down"

# Row 4: brand new data row appended to the sentiment log
$ws.Range("A4").Value = 42633.890567129631
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 26
$ws.Range("E4").Value = 22940
$ws.Range("F4").Value = 2661
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 88
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 16213
$ws.Range("L4").Value = 343
$ws.Range("M4").Value = 171
$ws.Range("N4").Value = 38
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = "Bag"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.86
$ws.Range("S4").Value = 0.026200000000000001
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = -2.66
$ws.Range("U4").Value = 15.05
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
